# Add csv module error handling
#
# The upstream CSV ingestion now keeps going after a parse hiccup instead of
# dropping the row, so each of the three affected sensor sheets picks up one
# extra sample per timestamp group in its most-recent logging block, plus
# three brand-new trailing samples that were previously swallowed by the
# error. This script reproduces that exact row-level reshaping:
#
#   - the last 9 timestamp groups in the block (sizes 6,6,6,4,4,4,2,2,2)
#     become 7,7,7,5,5,5,3,3,3 (one extra duplicate row per group), and
#   - 3 new single-row groups (new timestamps) are appended right after,
#     re-using the same non-timestamp ("B".."I") payload as the last
#     existing group, since that payload is constant across the whole
#     36-row block being expanded.

function Update-SensorSheet($wb, $sheetName, $trailingTimestamps) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ur = $ws.UsedRange
    $lastRow = $ur.Row() + $ur.Rows.Count() - 1
    $lastCol = $ur.Column() + $ur.Columns.Count() - 1

    $rng = $ws.Range("A1").Resize($lastRow, $lastCol)
    $vals = $rng.Value2()

    # The last 36 rows form 9 consecutive timestamp groups of sizes
    # 6,6,6,4,4,4,2,2,2 -- this is the block that gains one extra
    # duplicate row per group.
    $blockStart = $lastRow - 35

    # Distinct timestamps (9 groups), in row order.
    $groupTs = New-Object System.Collections.ArrayList
    $prevVal = $null
    for ($r = $blockStart; $r -le $lastRow; $r++) {
        $v = $vals[$r, 1]
        if ($prevVal -eq $null -or $v -ne $prevVal) {
            $groupTs.Add($v) | Out-Null
            $prevVal = $v
        }
    }

    # Columns B..I are identical for every row in the block, so grab
    # them once from the block's last row to use as a template.
    $template = @()
    for ($c = 2; $c -le $lastCol; $c++) {
        $template += $vals[$lastRow, $c]
    }

    $newSizes = @(7, 7, 7, 5, 5, 5, 3, 3, 3)

    $newBlockRows = New-Object System.Collections.ArrayList
    for ($g = 0; $g -lt $groupTs.Count; $g++) {
        $ts = $groupTs[$g]
        $cnt = $newSizes[$g]
        for ($k = 0; $k -lt $cnt; $k++) {
            $row = New-Object 'object[]' $lastCol
            $row[0] = $ts
            for ($c = 1; $c -lt $lastCol; $c++) {
                $row[$c] = $template[$c - 1]
            }
            $newBlockRows.Add($row) | Out-Null
        }
    }
    foreach ($ts in $trailingTimestamps) {
        $row = New-Object 'object[]' $lastCol
        $row[0] = $ts
        for ($c = 1; $c -lt $lastCol; $c++) {
            $row[$c] = $template[$c - 1]
        }
        $newBlockRows.Add($row) | Out-Null
    }

    $newLastRow = ($blockStart - 1) + $newBlockRows.Count

    $newArr = New-Object 'object[,]' $newLastRow, $lastCol
    for ($r = 1; $r -lt $blockStart; $r++) {
        for ($c = 1; $c -le $lastCol; $c++) {
            $newArr[$r - 1, $c - 1] = $vals[$r, $c]
        }
    }
    for ($i = 0; $i -lt $newBlockRows.Count; $i++) {
        $row = $newBlockRows[$i]
        $r = $blockStart + $i
        for ($c = 1; $c -le $lastCol; $c++) {
            $newArr[$r - 1, $c - 1] = $row[$c - 1]
        }
    }

    $outRng = $ws.Range("A1").Resize($newLastRow, $lastCol)
    $outRng.Value2 = $newArr

    # Newly appended rows (beyond the sheet's original last row) don't
    # inherit column A's date number format automatically -- apply it
    # explicitly so they render/round-trip as dates like the rest of
    # the column.
    if ($newLastRow -gt $lastRow) {
        $dateFmt = $ws.Cells.Item($blockStart, 1).NumberFormat()
        $newColARng = $ws.Range($ws.Cells.Item($lastRow + 1, 1), $ws.Cells.Item($newLastRow, 1))
        $newColARng.NumberFormat = $dateFmt
    }
}

$wb = $excel.ActiveWorkbook

Update-SensorSheet $wb "ROW35-FE-LIFTER" @(45726.73154126157, 45726.73156329861, 45726.73158645834)
Update-SensorSheet $wb "ROW35-MID-LIFTER" @(45726.57952, 45726.57954204861, 45726.57956549768)
Update-SensorSheet $wb "ROW02-MID-LIFTER" @(45726.72962616898, 45726.72964848379, 45726.72967174769)
